$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 7 with pull-sheet data (scaling + spiraling) ---
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "LOCAL"
$ws.Range("C7").Value = "543+00"
$ws.Range("D7").Value = "553+00"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "CABLE 1"
$ws.Range("G7").Value = "PK"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# E7 / G7 should pick up the same "centered, thin-bordered" look used by
# the rows above (style index shared with E2:E6 / G2:G6) rather than the
# plain style that row 7 started out with.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

# D7 is a brand-new cell so it naturally has no explicit style; F7 carried
# over row 7's old blank-cell style, so strip it back to Normal/default.
$ws.Range("F7").Style = "Normal"

$excel.CutCopyMode = 0

# --- Selection moves to F8 ---
$ws.Range("F8").Select()

# --- Window size/position change (zoom out of full screen) ---
$excel.Width = 13125
$excel.Height = 11550
$excel.Left = 4500
$excel.Top = 3180
